$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only these columns actually differ between row 8 and row 9 in the diff;
# swap just those so columns that already match (dates, comments, etc.)
# are left completely untouched (avoids COM re-typing text dates as
# numeric date serials on write-back).
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18, 26, 28)  # A,B,D,E,F,G,H,Q,R,Z,AB

foreach ($c in $cols) {
    $cell8 = $ws.Cells.Item(8, $c)
    $cell9 = $ws.Cells.Item(9, $c)

    $v8 = $cell8.Value2
    $v9 = $cell9.Value2

    $cell8.Value2 = $v9
    $cell9.Value2 = $v8
}
